$wb = $excel.ActiveWorkbook

# Update year-to-date (2024, column K) violent crime figures to reflect data through 2024-12-24.
# A handful of rows also touch columns D (2017) and I (2022) where historical corrections were made.

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 7769
$ws.Range("K3").Value = 8026
$ws.Range("D4").Value = 1955
$ws.Range("I4").Value = 1788
$ws.Range("K4").Value = 1688
$ws.Range("K5").Value = 575
$ws.Range("K6").Value = 8958
$ws.Range("D7").Value = 27755
$ws.Range("I7").Value = 25796
$ws.Range("K7").Value = 27016

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 493
$ws.Range("K3").Value = 533
$ws.Range("K6").Value = 594
$ws.Range("K7").Value = 1767

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 199
$ws.Range("K3").Value = 200
$ws.Range("K6").Value = 136
$ws.Range("K7").Value = 575

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 285
$ws.Range("K3").Value = 399
$ws.Range("K5").Value = 33
$ws.Range("K7").Value = 1133

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 150
$ws.Range("K7").Value = 449

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 253
$ws.Range("K5").Value = 34
$ws.Range("K7").Value = 889

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K3").Value = 34
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K6").Value = 195
$ws.Range("K7").Value = 807
$ws.Range("K8").Value = 1767
$ws.Range("K9").Value = 129
$ws.Range("K11").Value = 475
$ws.Range("K15").Value = 277
$ws.Range("K19").Value = 777
$ws.Range("K20").Value = 665
$ws.Range("K25").Value = 130
$ws.Range("K29").Value = 1493
$ws.Range("K30").Value = 111
$ws.Range("K31").Value = 324
$ws.Range("K33").Value = 1133
$ws.Range("K36").Value = 350
$ws.Range("K37").Value = 889
$ws.Range("K42").Value = 1001
$ws.Range("K43").Value = 225
$ws.Range("K46").Value = 55
$ws.Range("K47").Value = 185
$ws.Range("K48").Value = 337
$ws.Range("K51").Value = 349
$ws.Range("K54").Value = 528
$ws.Range("K55").Value = 296
$ws.Range("D63").Value = 353
$ws.Range("I63").Value = 235
$ws.Range("K66").Value = 80
$ws.Range("K67").Value = 1048
$ws.Range("K76").Value = 374
$ws.Range("K79").Value = 660
$ws.Range("K83").Value = 575
$ws.Range("K85").Value = 1243
$ws.Range("K86").Value = 164
$ws.Range("K88").Value = 287
$ws.Range("K89").Value = 406
$ws.Range("K90").Value = 258
$ws.Range("K91").Value = 324
$ws.Range("K95").Value = 449
$ws.Range("K96").Value = 293
$ws.Range("D101").Value = 27755
$ws.Range("I101").Value = 25796
$ws.Range("K101").Value = 27016

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K3").Value = 79
$ws.Range("K6").Value = 128
$ws.Range("K7").Value = 324

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 381
$ws.Range("K6").Value = 297
$ws.Range("K7").Value = 1048

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 286
$ws.Range("K7").Value = 528

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 419
$ws.Range("K3").Value = 525
$ws.Range("K4").Value = 69
$ws.Range("K6").Value = 443
$ws.Range("K7").Value = 1493

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K2").Value = 52
$ws.Range("K7").Value = 337

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 227
$ws.Range("K3").Value = 232
$ws.Range("K5").Value = 25
$ws.Range("K6").Value = 259
$ws.Range("K7").Value = 777

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K2").Value = 83
$ws.Range("K3").Value = 74
$ws.Range("K7").Value = 374

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 76
$ws.Range("K7").Value = 195

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K4").Value = 44
$ws.Range("K7").Value = 1001

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K2").Value = 87
$ws.Range("K7").Value = 296

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K2").Value = 89
$ws.Range("K3").Value = 62
$ws.Range("K6").Value = 117
$ws.Range("K7").Value = 293

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 149
$ws.Range("K7").Value = 324

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K3").Value = 209
$ws.Range("K7").Value = 660

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K6").Value = 193
$ws.Range("K7").Value = 665

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 350

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K6").Value = 229
$ws.Range("K7").Value = 807

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 130

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K4").Value = 19
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 185

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 277

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 165
$ws.Range("K7").Value = 475

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K3").Value = 47
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 90
$ws.Range("K7").Value = 287

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K6").Value = 122
$ws.Range("K7").Value = 406

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 164

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K4").Value = 17
$ws.Range("K7").Value = 258

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K4").Value = 39
$ws.Range("K7").Value = 349

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K4").Value = 30
$ws.Range("K7").Value = 225

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 411
$ws.Range("K3").Value = 430
$ws.Range("K7").Value = 1243
